# Update "想去人数" (column F) counts on the "展览" and "全部类型" sheets
# to reflect freshly re-generated data (gh-pages output regenerated).

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item(1)  # 展览
$sheetAllTypes   = $wb.Worksheets.Item(4)  # 全部类型

# Row -> new F value for the "展览" sheet
$exhibitionUpdates = @{
    3  = 12251
    4  = 4498
    6  = 70
    7  = 33
    9  = 2616
    10 = 1137
    13 = 5355
    14 = 67
    15 = 218
    16 = 566
    17 = 11505
    18 = 11617
    23 = 58
}

foreach ($row in $exhibitionUpdates.Keys) {
    $sheetExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

# Row -> new F value for the "全部类型" sheet
$allTypesUpdates = @{
    3  = 12251
    4  = 4498
    6  = 70
    7  = 33
    9  = 2617
    11 = 1137
    14 = 5355
    15 = 67
    16 = 218
    17 = 566
    18 = 11505
    19 = 11617
    24 = 58
}

foreach ($row in $allTypesUpdates.Keys) {
    $sheetAllTypes.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
